$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.57880000000002
$ws.Range("B4").Value = 5.296699999999996
$ws.Range("C4").Value = -11.42949999999999
$ws.Range("E4").Value = 13.974

$ws.Range("B5").Value = 5.225399999999997

$ws.Range("A7").Value = -21.5406

$ws.Range("B8").Value = 5.047899999999999

$ws.Range("C9").Value = -11.74570000000001

$ws.Range("E12").Value = 11.7043

$ws.Range("A16").Value = -21.52600000000002
$ws.Range("B16").Value = 4.838000000000004

$ws.Range("C18").Value = -14.34880000000001

$ws.Range("E20").Value = 13.14199999999999
